$d = $word.ActiveDocument

function Split-Biblio($anchor, $replacement) {
    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Replacement.ClearFormatting()
    $find.Execute($anchor, $false, $false, $false, $false, $false, $true, 1, $false, $replacement, 2) | Out-Null
}

Split-Biblio "2005.JAFFE"        "2005.^l^lJAFFE"
Split-Biblio "2002.GITMAN"       "2002.^l^lGITMAN"
Split-Biblio "2017.ASSAF"        "2017.^l^lASSAF"
Split-Biblio "2014MORANTE"       "2014^l^lMORANTE"
Split-Biblio "2009.NEWNAN"       "2009.^l^lNEWNAN"
Split-Biblio "2000.HOJI"         "2000.^l^lHOJI"
Split-Biblio "2019.SANVICENTE"   "2019.^l^lSANVICENTE"

Write-Host "done"
